# Update the two YouTube link rows and turn the second one (A2) into a
# real hyperlink, so comment/reply order can be keyed off it.
#
#   A1: https://www.youtube.com/watch?v=0UKwpJUUDlM  -> https://www.youtube.com/shorts/xspWfbOSsms
#   A2: https://www.youtube.com/shorts/qmQ_lzNtQco   -> https://www.youtube.com/watch?v=eb6cJMSZuWo (+ hyperlink)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "https://www.youtube.com/shorts/xspWfbOSsms"
$ws.Range("A2").Value = "https://www.youtube.com/watch?v=eb6cJMSZuWo"

# Adding the hyperlink drives Excel's own "Hyperlink" cell style (new font +
# xf + cellStyle entries in styles.xml) onto A2, matching a manual
# Insert > Hyperlink on that cell.
$ws.Hyperlinks.Add($ws.Range("A2"), "https://www.youtube.com/watch?v=eb6cJMSZuWo")
